# C5-PowerPoint.pptx edit:
#  1. Slide 6's table switches from the custom "Table_0" style to the
#     built-in table style {F90B60FD-4CB7-4C4D-9432-BE31291BA696}.
#  2. The deck's theme (serialized to ppt/theme/theme1.xml, shared by the
#     slide master) is recolored from the "Integral" palette to the
#     stock "Office Theme" palette (font scheme / format scheme are
#     already identical between the two themes, only clrScheme differs).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{F90B60FD-4CB7-4C4D-9432-BE31291BA696}")

# --- 2) Theme colors: Integral -> Office Theme -----------------------------
$slideMaster = $p.Slides.Item(1).Design.SlideMaster
$colorScheme = $slideMaster.Theme.ThemeColorScheme

# Colors(1..12) == dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$colorScheme.Item(1).RGB  = 0        # dk1      000000
$colorScheme.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388  # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501  # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407    # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308 # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456  # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797 # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477  # folHlink 954F72
